$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 428179.38
$ws.Range("J17").Value = 428179.38
$ws.Range("L17").Value = 1284538.14
$ws.Range("N17").Value = -1284874.14
$ws.Range("H28").Value = 952.0625
$ws.Range("I28").Value = 1374.5555
$ws.Range("J28").Value = 408.85715
$ws.Range("K28").Value = 1374.5555
$ws.Range("L28").Value = 408.85715
$ws.Range("M28").Value = -889.5554999999999
$ws.Range("N28").Value = -1378.85715
$ws.Range("H40").Value = 2207.8948
$ws.Range("I40").Value = 1712.5
$ws.Range("J40").Value = 2340
$ws.Range("K40").Value = 1712.5
$ws.Range("L40").Value = 2340
$ws.Range("M40").Value = -1537.5
$ws.Range("N40").Value = -2690
$ws.Range("H62").Value = 95613.336
$ws.Range("I62").Value = 111485.4
$ws.Range("J62").Value = 16253
$ws.Range("K62").Value = 111485.4
$ws.Range("L62").Value = 16253
$ws.Range("M62").Value = -110861.4
$ws.Range("N62").Value = -17501
$ws.Range("H65").Value = 95613.336
$ws.Range("I65").Value = 111485.4
$ws.Range("J65").Value = 16253
$ws.Range("K65").Value = 557427
$ws.Range("L65").Value = 81265
$ws.Range("M65").Value = -554307
$ws.Range("N65").Value = -87505
$ws.Range("H106").Value = 2738.9333
$ws.Range("I106").Value = 2560.4
$ws.Range("J106").Value = 3096
$ws.Range("K106").Value = 2560.4
$ws.Range("L106").Value = 3096
$ws.Range("M106").Value = -1929.4
$ws.Range("N106").Value = -4358
$ws.Range("H132").Value = 4260.1
$ws.Range("I132").Value = 4961.875
$ws.Range("J132").Value = 1453
$ws.Range("K132").Value = 14885.625
$ws.Range("L132").Value = 4359
$ws.Range("M132").Value = -12355.625
$ws.Range("N132").Value = -9419
$ws.Range("H137").Value = 12821984
$ws.Range("I137").Value = 1565.7858
$ws.Range("K137").Value = 4697.357400000001
$ws.Range("M137").Value = -2147.357400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 896.7917
$ws.Range("I74").Value = 959.36365
$ws.Range("J74").Value = 843.8461
$ws.Range("K74").Value = 959.36365
$ws.Range("L74").Value = 843.8461
$ws.Range("M74").Value = -85.36365000000001
$ws.Range("N74").Value = -2591.8461
$ws.Range("H77").Value = 896.7917
$ws.Range("I77").Value = 959.36365
$ws.Range("J77").Value = 843.8461
$ws.Range("K77").Value = 4796.81825
$ws.Range("L77").Value = 4219.2305
$ws.Range("M77").Value = -428.8182500000003
$ws.Range("N77").Value = -12955.2305
$ws.Range("H97").Value = 945.13336
$ws.Range("I97").Value = 782.8461
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 782.8461
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -286.8461
$ws.Range("N97").Value = -2992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 29225.666
$ws.Range("J87").Value = 29225.666
$ws.Range("L87").Value = 29225.666
$ws.Range("N87").Value = -31721.666
$ws.Range("H90").Value = 29225.666
$ws.Range("J90").Value = 29225.666
$ws.Range("L90").Value = 87676.99800000001
$ws.Range("N90").Value = -100156.998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2248.2856
$ws.Range("I31").Value = 1866.4117
$ws.Range("J31").Value = 2608.9443
$ws.Range("K31").Value = 1866.4117
$ws.Range("L31").Value = 2608.9443
$ws.Range("M31").Value = -1571.4117
$ws.Range("N31").Value = -3198.9443
$ws.Range("H34").Value = 2248.2856
$ws.Range("I34").Value = 1866.4117
$ws.Range("J34").Value = 2608.9443
$ws.Range("K34").Value = 1866.4117
$ws.Range("L34").Value = 2608.9443
$ws.Range("M34").Value = -1664.4117
$ws.Range("N34").Value = -3012.9443
$ws.Range("H140").Value = 46252.727
$ws.Range("J140").Value = 46252.727
$ws.Range("L140").Value = 46252.727
$ws.Range("N140").Value = -56612.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 5560
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5560
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 16680
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -17536
$ws.Range("H91").Value = 5560
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5560
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 16680
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -19644
$ws.Range("H113").Value = 400.44
$ws.Range("I113").Value = 369.5
$ws.Range("J113").Value = 406.33334
$ws.Range("K113").Value = 1108.5
$ws.Range("L113").Value = 1219.00002
$ws.Range("M113").Value = 1061.5
$ws.Range("N113").Value = -5559.000019999999
$ws.Range("H131").Value = 15636.057
$ws.Range("J131").Value = 1614.3103
$ws.Range("L131").Value = 4842.9309
$ws.Range("N131").Value = -14922.9309
$ws.Range("H132").Value = 762.5
$ws.Range("I132").Value = 723.0769
$ws.Range("J132").Value = 933.3333
$ws.Range("K132").Value = 6507.6921
$ws.Range("L132").Value = 8399.9997
$ws.Range("M132").Value = -3977.6921
$ws.Range("N132").Value = -13459.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1200.7646
$ws.Range("I102").Value = 922.3570999999999
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 922.3570999999999
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = 699.6429000000001
$ws.Range("N102").Value = -5744
$ws.Range("H132").Value = 3663.5789
$ws.Range("I132").Value = 3126.3333
$ws.Range("J132").Value = 4584.5713
$ws.Range("K132").Value = 9378.999899999999
$ws.Range("L132").Value = 13753.7139
$ws.Range("M132").Value = -6848.999899999999
$ws.Range("N132").Value = -18813.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 478.81818
$ws.Range("I46").Value = 392
$ws.Range("J46").Value = 528.4286
$ws.Range("K46").Value = 392
$ws.Range("L46").Value = 528.4286
$ws.Range("M46").Value = -204
$ws.Range("N46").Value = -904.4286
$ws.Range("H132").Value = 3208064.2
$ws.Range("I132").Value = 4169143.8
$ws.Range("K132").Value = 12507431.4
$ws.Range("M132").Value = -12504901.4
$ws.Range("H136").Value = 4123.375
$ws.Range("J136").Value = 5490
$ws.Range("L136").Value = 16470
$ws.Range("N136").Value = -21570

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 7000
$ws.Range("J53").Value = 7000
$ws.Range("L53").Value = 7000
$ws.Range("N53").Value = -8214
$ws.Range("H132").Value = 2229.7334
$ws.Range("I132").Value = 1606
$ws.Range("J132").Value = 2775.5
$ws.Range("K132").Value = 4818
$ws.Range("L132").Value = 8326.5
$ws.Range("M132").Value = -2288
$ws.Range("N132").Value = -13386.5
$ws.Range("H136").Value = 2237.1555
$ws.Range("I136").Value = 2382.577
$ws.Range("J136").Value = 2038.1578
$ws.Range("K136").Value = 7147.731000000001
$ws.Range("L136").Value = 6114.4734
$ws.Range("M136").Value = -4597.731000000001
$ws.Range("N136").Value = -11214.4734
